$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell updates as literal text, matching the source data format
# (prices/links/percentages/coin names are stored as plain text strings).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.256.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.428.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +7.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.41"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.16"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.801.66"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.421.73"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +7.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.848"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.13"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.158.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.82"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +7.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.69"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.33%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.57"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.63%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.80"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.84"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +14.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.75"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.57"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.86"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0773"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.01"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +18.93%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.09"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.90"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0305"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.27"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.979.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +11.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.96"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.41"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +30.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.70"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +10.89%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.15%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.188"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.669.65"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.26%  "
